# Add the new "Anosim" sheet right after "adon.results"
$wb = $excel.ActiveWorkbook
$anchor = $wb.Worksheets.Item("adon.results")
$anosim = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $anchor)
$anosim.Name = "Anosim"

# Populate new shared strings in the order they appear in the target file
# (100=anosim_r, 101=R_value, 102=p_value, 103=ANOSIM desc)
$anosim.Range("A2").Value = "anosim_r"
$anosim.Range("B3").Value = "R_value"
$anosim.Range("C3").Value = "p_value"
$anosim.Range("A1").Value = "ANOSIM for region by year, 10000 permutations Bray-Curtis"
$anosim.Range("A3").Value = "Year"

# A2 uses the same "Lucida Console, 10pt, blue" style already used by the
# "> modtab_cal" / "> modtab_cyclo" / "> modtab_clad" header cells elsewhere
# in the workbook - copy that formatting over instead of re-creating it.
$wb.Worksheets.Item("totalzoop_posthoc").Range("A1").Copy()
$anosim.Range("A2").PasteSpecial(-4122)

$yearData = @(
    @(2014, 0.28367199999999998, 0.00019998),
    @(2015, 0.32816679999999998, 0.00009999),
    @(2016, 0.28499099999999999, 0.00009999),
    @(2017, 0.1723085,           0.00029997),
    @(2018, 0.25141419999999998, 0.00009999),
    @(2019, 0.54593349999999996, 0.00009999)
)
$r = 4
foreach ($row in $yearData) {
    $anosim.Cells.Item($r, 1).Value = $row[0]
    $anosim.Cells.Item($r, 2).Value = $row[1]
    $anosim.Cells.Item($r, 3).Value = $row[2]
    $r++
}

$anosim.Columns.Item(3).ColumnWidth = 13.1666666666667
$anosim.Range("B16").Select()

# Prepend a title row and the R code row to "adon.results"
$adon = $wb.Worksheets.Item("adon.results")
$adon.Rows("1:2").Insert()
$adon.Range("A1").Value = "PERMANOVA of [TaxonName zoop cpue] by region, year, sampleperiod with interactions"
$adon.Range("A2").Value = 'adon.results <- adonis2(genw2[c(9:70)]~genw2$Region*genw2$SamplePeriod+genw2$Year*genw2$SamplePeriod+genw2$Region*genw2$Year, strata=genw2$StationCode,method=bray,perm=999)'
$adon.Columns.Item(1).ColumnWidth = 26.7
$adon.Columns.Item(9).ColumnWidth = 8.3
$adon.Range("I8").Select()

$adon.Activate()
